# Add a new data row (row 6) to the Artfynd sheet, mirroring the structure
# of the existing rows but for a new observation (Tjäder / Tetrao urogallus).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Numeric cells ---
$ws.Range("A6").Value = 131258730
$ws.Range("B6").Value = 57073
$ws.Range("E6").Value = 100138
$ws.Range("Q6").Value = 540241
$ws.Range("R6").Value = 6737560
$ws.Range("S6").Value = 10

# --- Text cells ---
$ws.Range("D6").Value = "LC"
$ws.Range("F6").Value = "Tjäder"
$ws.Range("G6").Value = "Tetrao urogallus"
$ws.Range("H6").Value = "Linnaeus, 1758"
# "Antal" looks numeric but must be stored as text, like the source data.
$ws.Range("I6").Value = "'2"
$ws.Range("M6").Value = "övernattning"
$ws.Range("P6").Value = "Skallberget, Skallberget, Dlr"
$ws.Range("T6").Value = "Dalarna"
$ws.Range("U6").Value = "Falun"
$ws.Range("V6").Value = "Dalarna"
$ws.Range("W6").Value = "Svärdsjö"
# Dates/times are stored as plain text in this sheet, not as Excel date serials.
$ws.Range("Y6").Value = "'2026-02-22"
$ws.Range("Z6").Value = "13:20"
$ws.Range("AA6").Value = "'2026-02-22"
$ws.Range("AB6").Value = "13:20"
$ws.Range("AW6").Value = "Göran Ehn"
$ws.Range("AX6").Value = "Göran Ehn"

# --- Boolean cells ---
$ws.Range("AD6").Value = $false
$ws.Range("AE6").Value = $false
$ws.Range("AG6").Value = $false
